# Generate Report for Handback
#
# Row 7 (fdec13da-2ffe-43d8-a64a-daf16cbeac71) on both the "zh-cn" and
# "de-de" sheets had its handback evaluated: Excel found that the handback
# file supplied is out of date versus the latest source, so:
#   - "Latest Target File" (I7) now links to the target .md file
#   - "Latest Handback File" (J7) is filled in with the xlf file name
#   - "Latest Handback DateTime" (K7) gets the detection timestamp
#   - "Error Detail" (P7) explains that the handback file is stale

$wb = $excel.ActiveWorkbook

$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4cd98338e337d8976c80311f3dc741b05d296c0/e2e/fdec13da-2ffe-43d8-a64a-daf16cbeac71.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7ca877445c78732b3ccf72653ca702f942d9ce6d/e2e/fdec13da-2ffe-43d8-a64a-daf16cbeac71.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4cd98338e337d8976c80311f3dc741b05d296c0/e2e/fdec13da-2ffe-43d8-a64a-daf16cbeac71.md."

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("J7").Value = "fdec13da-2ffe-43d8-a64a-daf16cbeac71.c2c44bc92d2b9d4d546c27547f7757ffe6045d56.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-09-01 00:59:11"
$wsZhCn.Range("P7").Value = $errorDetail

$wsZhCn.Range("I7").Value = "fdec13da-2ffe-43d8-a64a-daf16cbeac71.md"
$wsZhCn.Range("I7").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $latestMdUrl, "", "", "fdec13da-2ffe-43d8-a64a-daf16cbeac71.md")

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("J7").Value = "fdec13da-2ffe-43d8-a64a-daf16cbeac71.c2c44bc92d2b9d4d546c27547f7757ffe6045d56.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-09-01 00:59:19"
$wsDeDe.Range("P7").Value = $errorDetail

$wsDeDe.Range("I7").Value = "fdec13da-2ffe-43d8-a64a-daf16cbeac71.md"
$wsDeDe.Range("I7").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $latestMdUrl, "", "", "fdec13da-2ffe-43d8-a64a-daf16cbeac71.md")
